# Data refresh: update "想去人数" (column F) counts across sheets
# as generated at commit 456a3b4 (gh-pages output update).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F5").Value = 8851
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 7095
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 5389
$ws.Range("F12").Value = 6094
$ws.Range("F13").Value = 1095
$ws.Range("F14").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("F20").Value = 144
$ws.Range("F21").Value = 202
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 9968
$ws.Range("F25").Value = 1881
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 80
$ws.Range("F31").Value = 86
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 1029
$ws.Range("F34").Value = 0
$ws.Range("F36").Value = 1358
$ws.Range("F38").Value = 0
$ws.Range("F40").Value = 1194
$ws.Range("F42").Value = 101
$ws.Range("F43").Value = 168
$ws.Range("F46").Value = 973
$ws.Range("F47").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("F49").Value = 0

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("F18").Value = 904
$ws.Range("F21").Value = 0

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("F9").Value = 7096
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 5389
$ws.Range("F15").Value = 6094
$ws.Range("F16").Value = 1095
$ws.Range("F18").Value = 399
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 328
$ws.Range("F22").Value = 270
$ws.Range("F23").Value = 0
$ws.Range("F27").Value = 190
$ws.Range("F28").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 1769
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 86
$ws.Range("F37").Value = 2031
$ws.Range("F39").Value = 1358
$ws.Range("F41").Value = 0
$ws.Range("F42").Value = 643
$ws.Range("F43").Value = 0
$ws.Range("F44").Value = 168
$ws.Range("F46").Value = 1072
$ws.Range("F48").Value = 1353
$ws.Range("F49").Value = 0
$ws.Range("F50").Value = 0
